$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (BGR ints as used by Excel COM Interior.Color)
$green  = 5287936    # FF00B050
$yellow = 65535       # FFFFFF00
$purple = 10498160    # FF7030A0

# Row 15: ESTADO PROGRESS (yellow) -> DONE (green)
$ws.Range("F15").Interior.Color = $green
$ws.Range("F15").Value = "DONE"

# Row 16: ESTADO FLUTTER (blue) -> DONE (green)
$ws.Range("F16").Interior.Color = $green
$ws.Range("F16").Value = "DONE"

# Row 18: ESTADO DONE stays DONE, but fill changes yellow -> green
$ws.Range("F18").Interior.Color = $green
$ws.Range("F18").Value = "DONE"

# Row 19: ESTADO FLUTTER (blue) -> PROGRESS (yellow)
$ws.Range("F19").Interior.Color = $yellow
$ws.Range("F19").Value = "PROGRESS"

# Row 31: TIPO POST (orange) -> "-" (new purple fill/style)
$ws.Range("D31").Interior.Color = $purple
$ws.Range("D31").Value = "-"

# Move active selection to F19, matching the saved cursor position
$ws.Range("F19").Select()
